# Weekly update: insert two new price records for Pimiento (Zafiro rojo / Zafiro verde)
# at Macroferia Regional de Talca, pushing the existing historical rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 214 (this shifts rows 214..312 down to 216..314)
$ws.Rows.Item(214).Insert()
$ws.Rows.Item(214).Insert()

# New row 214: Zafiro rojo
$ws.Cells.Item(214, 1).Value = 5
$ws.Cells.Item(214, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(214, 3).Value = "Maule"
$ws.Cells.Item(214, 4).Value = 44455
$ws.Cells.Item(214, 5).Value = 7
$ws.Cells.Item(214, 6).Value = 100112002
$ws.Cells.Item(214, 7).Value = "Pimiento"
$ws.Cells.Item(214, 8).Value = "Zafiro rojo"
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 200
$ws.Cells.Item(214, 11).Value = 38000
$ws.Cells.Item(214, 12).Value = 38000
$ws.Cells.Item(214, 13).Value = 38000
$ws.Cells.Item(214, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(214, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(214, 16).Value = 2533
$ws.Cells.Item(214, 17).Value = 15
$ws.Cells.Item(214, 18).Value = "Hortaliza"

# New row 215: Zafiro verde
$ws.Cells.Item(215, 1).Value = 5
$ws.Cells.Item(215, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(215, 3).Value = "Maule"
$ws.Cells.Item(215, 4).Value = 44455
$ws.Cells.Item(215, 5).Value = 7
$ws.Cells.Item(215, 6).Value = 100112002
$ws.Cells.Item(215, 7).Value = "Pimiento"
$ws.Cells.Item(215, 8).Value = "Zafiro verde"
$ws.Cells.Item(215, 9).Value = "Primera"
$ws.Cells.Item(215, 10).Value = 200
$ws.Cells.Item(215, 11).Value = 35000
$ws.Cells.Item(215, 12).Value = 35000
$ws.Cells.Item(215, 13).Value = 35000
$ws.Cells.Item(215, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(215, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(215, 16).Value = 2333
$ws.Cells.Item(215, 17).Value = 15
$ws.Cells.Item(215, 18).Value = "Hortaliza"
